$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1225
$ws.Range("I18").Value = 1314.2858
$ws.Range("J18").Value = 600
$ws.Range("K18").Value = 1314.2858
$ws.Range("L18").Value = 600
$ws.Range("M18").Value = -1030.2858
$ws.Range("N18").Value = -1168

$ws.Range("H40").Value = 1064.9783
$ws.Range("I40").Value = 1055.5588
$ws.Range("J40").Value = 1091.6666
$ws.Range("K40").Value = 1055.5588
$ws.Range("L40").Value = 1091.6666
$ws.Range("M40").Value = -880.5588
$ws.Range("N40").Value = -1441.6666

$ws.Range("H64").Value = 3333.2666
$ws.Range("I64").Value = 3205.077
$ws.Range("J64").Value = 4166.5
$ws.Range("K64").Value = 3205.077
$ws.Range("L64").Value = 4166.5
$ws.Range("M64").Value = -2957.077
$ws.Range("N64").Value = -4662.5

$ws.Range("H67").Value = 3333.2666
$ws.Range("I67").Value = 3205.077
$ws.Range("J67").Value = 4166.5
$ws.Range("K67").Value = 3205.077
$ws.Range("L67").Value = 4166.5
$ws.Range("M67").Value = -2347.077
$ws.Range("N67").Value = -5882.5

$ws.Range("H74").Value = 3452.5652
$ws.Range("I74").Value = 3349.9375
$ws.Range("J74").Value = 3687.1428
$ws.Range("K74").Value = 3349.9375
$ws.Range("L74").Value = 3687.1428
$ws.Range("M74").Value = -2413.9375
$ws.Range("N74").Value = -5559.1428

$ws.Range("H77").Value = 3452.5652
$ws.Range("I77").Value = 3349.9375
$ws.Range("J77").Value = 3687.1428
$ws.Range("K77").Value = 16749.6875
$ws.Range("L77").Value = 18435.714
$ws.Range("M77").Value = -12069.6875
$ws.Range("N77").Value = -27795.714

$ws.Range("H138").Value = 2387.4478
$ws.Range("I138").Value = 1143.65
$ws.Range("J138").Value = 4230.1113
$ws.Range("K138").Value = 3430.95
$ws.Range("L138").Value = 12690.3339
$ws.Range("M138").Value = 1709.05
$ws.Range("N138").Value = -22970.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = $null
$ws.Range("N5").Value = -724

$ws.Range("H32").Value = 5389.017
$ws.Range("I32").Value = 4132.0586
$ws.Range("K32").Value = 4132.0586
$ws.Range("M32").Value = -3845.0586

$ws.Range("H39").Value = 33009
$ws.Range("I39").Value = 3000
$ws.Range("J39").Value = 63018
$ws.Range("K39").Value = 3000
$ws.Range("L39").Value = 63018
$ws.Range("M39").Value = -2480
$ws.Range("N39").Value = -64058

$ws.Range("H63").Value = 4180
$ws.Range("I63").Value = 4176
$ws.Range("J63").Value = 4200
$ws.Range("K63").Value = 4176
$ws.Range("L63").Value = 4200
$ws.Range("M63").Value = -3490
$ws.Range("N63").Value = -5572

$ws.Range("H66").Value = 4180
$ws.Range("I66").Value = 4176
$ws.Range("J66").Value = 4200
$ws.Range("K66").Value = 20880
$ws.Range("L66").Value = 21000
$ws.Range("M66").Value = -17448
$ws.Range("N66").Value = -27864

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = -730

$ws.Range("H15").Value = 21647.857
$ws.Range("J15").Value = 21647.857
$ws.Range("L15").Value = 21647.857
$ws.Range("N15").Value = -22101.857

$ws.Range("H35").Value = 15000
$ws.Range("J35").Value = 15000
$ws.Range("L35").Value = 15000
$ws.Range("N35").Value = -15620

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 940
$ws.Range("I22").Value = 985.38464
$ws.Range("K22").Value = 985.38464
$ws.Range("M22").Value = -635.38464

$ws.Range("H31").Value = 31084
$ws.Range("I31").Value = 51285.8
$ws.Range("J31").Value = 4148.2666
$ws.Range("K31").Value = 51285.8
$ws.Range("L31").Value = 4148.2666
$ws.Range("M31").Value = -50990.8
$ws.Range("N31").Value = -4738.2666

$ws.Range("H34").Value = 31084
$ws.Range("I34").Value = 51285.8
$ws.Range("J34").Value = 4148.2666
$ws.Range("K34").Value = 51285.8
$ws.Range("L34").Value = 4148.2666
$ws.Range("M34").Value = -51083.8
$ws.Range("N34").Value = -4552.2666

$ws.Range("H41").Value = 18931.334
$ws.Range("I41").Value = 5014.25
$ws.Range("J41").Value = 30065
$ws.Range("K41").Value = 5014.25
$ws.Range("L41").Value = 30065
$ws.Range("M41").Value = -4586.25
$ws.Range("N41").Value = -30921

$ws.Range("H50").Value = 20186.154
$ws.Range("I50").Value = 500
$ws.Range("J50").Value = 26092
$ws.Range("K50").Value = 500
$ws.Range("L50").Value = 26092
$ws.Range("M50").Value = 125
$ws.Range("N50").Value = -27342

$ws.Range("H51").Value = 25915.834
$ws.Range("I51").Value = 15000
$ws.Range("J51").Value = 28099
$ws.Range("K51").Value = 15000
$ws.Range("L51").Value = 28099
$ws.Range("M51").Value = -14264
$ws.Range("N51").Value = -29571

$ws.Range("H59").Value = 41224
$ws.Range("J59").Value = 41224
$ws.Range("L59").Value = 41224
$ws.Range("N59").Value = -43514

$ws.Range("H60").Value = 23266.092
$ws.Range("I60").Value = 12500
$ws.Range("J60").Value = 25658.555
$ws.Range("K60").Value = 12500
$ws.Range("L60").Value = 25658.555
$ws.Range("M60").Value = -11989
$ws.Range("N60").Value = -26680.555

$ws.Range("H61").Value = 25915.834
$ws.Range("I61").Value = 15000
$ws.Range("J61").Value = 28099
$ws.Range("K61").Value = 15000
$ws.Range("L61").Value = 28099
$ws.Range("M61").Value = -14652
$ws.Range("N61").Value = -28795

$ws.Range("H62").Value = 2774.7334
$ws.Range("I62").Value = 2721.5
$ws.Range("J62").Value = 2881.2
$ws.Range("K62").Value = 2721.5
$ws.Range("L62").Value = 2881.2
$ws.Range("M62").Value = -2097.5
$ws.Range("N62").Value = -4129.2

$ws.Range("H65").Value = 2774.7334
$ws.Range("I65").Value = 2721.5
$ws.Range("J65").Value = 2881.2
$ws.Range("K65").Value = 13607.5
$ws.Range("L65").Value = 14406
$ws.Range("M65").Value = -10487.5
$ws.Range("N65").Value = -20646

$ws.Range("H68").Value = 32368.525
$ws.Range("J68").Value = 32368.525
$ws.Range("L68").Value = 32368.525
$ws.Range("N68").Value = -33866.525

$ws.Range("H71").Value = 32368.525
$ws.Range("J71").Value = 32368.525
$ws.Range("L71").Value = 97105.57500000001
$ws.Range("N71").Value = -104593.575

$ws.Range("H74").Value = 36691.2
$ws.Range("J74").Value = 36691.2
$ws.Range("L74").Value = 36691.2
$ws.Range("N74").Value = -38439.2

$ws.Range("H77").Value = 36691.2
$ws.Range("J77").Value = 36691.2
$ws.Range("L77").Value = 110073.6
$ws.Range("N77").Value = -118809.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1063.6296
$ws.Range("I113").Value = 1818.6364
$ws.Range("J113").Value = 544.5625
$ws.Range("K113").Value = 5455.9092
$ws.Range("L113").Value = 1633.6875
$ws.Range("M113").Value = -3285.9092
$ws.Range("N113").Value = -5973.6875

$ws.Range("H133").Value = 6090.696
$ws.Range("I133").Value = 5794.875
$ws.Range("J133").Value = 6248.467
$ws.Range("K133").Value = 17384.625
$ws.Range("L133").Value = 18745.401
$ws.Range("M133").Value = -12324.625
$ws.Range("N133").Value = -28865.401

$ws.Range("H134").Value = 2204.8635
$ws.Range("I134").Value = 1028.1666
$ws.Range("J134").Value = 7500
$ws.Range("K134").Value = 3084.4998
$ws.Range("L134").Value = 22500
$ws.Range("M134").Value = 1985.5002
$ws.Range("N134").Value = -32640

$ws.Range("H137").Value = 2809.7368
$ws.Range("I137").Value = 1096.8
$ws.Range("J137").Value = 9233.25
$ws.Range("K137").Value = 3290.4
$ws.Range("L137").Value = 27699.75
$ws.Range("M137").Value = 1809.6
$ws.Range("N137").Value = -37899.75

$ws.Range("H139").Value = 1281.1111
$ws.Range("I139").Value = 816.25
$ws.Range("K139").Value = 2448.75
$ws.Range("M139").Value = 2691.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 500
$ws.Range("I46").Value = 500
$ws.Range("K46").Value = 500
$ws.Range("M46").Value = -312

$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -48180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 30000
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314
